# Add files via upload
# Recreates the "Somatoria" worksheet: a small two-column table (with bold
# header row) and the page setup used when it was re-saved by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (bold) --------------------------------------------------
$ws.Range("A1").Value = "Coluna1"
$ws.Range("B1").Value = "Couluna2"
$ws.Range("A1:B1").Font.Bold = $true

# --- Data rows (literal numbers entered as "=n" formulas, as in source) -
$ws.Range("A2").Formula = "=10"
$ws.Range("B2").Formula = "=5"

$ws.Range("A3").Formula = "=20"
$ws.Range("B3").Formula = "=2"

$ws.Range("A4").Formula = "=8"
$ws.Range("B4").Formula = "=3"

$ws.Range("A5").Formula = "=9"
$ws.Range("B5").Formula = "=50"

# --- Page setup used on export ------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Select the populated range, matching the saved view ---------------
[void]$ws.Range("A1:B5").Select()
